$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COVID19_TIMESERIESDATA")

# New header row (row 901) replicates the recurring block-header style/labels
# used throughout the sheet (e.g. the block starting at row 865), giving the
# same bold/bordered "States/UT | Active Cases | ... | Date" header.
$ws.Range("A865:H865").Copy($ws.Range("A901:H901"))

# New day's data (12-10-2020), one row per State/UT, appended right after
# the new header row, continuing the existing table layout.
$covidData = @(
    @("Andaman and Nicobar Islands", 186, -7, 3764, 20, 55, 0, "12-10-2020"),
    @("Andhra Pradesh", 46295, -329, 703208, 5509, 6224, 30, "12-10-2020"),
    @("Arunachal Pradesh", 2891, -49, 9232, 197, 24, 1, "12-10-2020"),
    @("Assam", 28385, -836, 164582, 1227, 816, 5, "12-10-2020"),
    @("Bihar", 11044, -121, 184541, 1151, 946, 2, "12-10-2020"),
    @("Chandigarh", 1184, -45, 11787, 125, 191, 1, "12-10-2020"),
    @("Chhattisgarh", 27348, -21, 113771, 2117, 1253, 18, "12-10-2020"),
    @("Dadra and Nagar Haveli and Daman and Diu", 102, 2, 3059, 9, 2, 0, "12-10-2020"),
    @("Delhi", 21701, -306, 281869, 3057, 5769, 29, "12-10-2020"),
    @("Goa", 4656, -2, 33203, 426, 507, 8, "12-10-2020"),
    @("Gujarat", 15695, -241, 132173, 1413, 3566, 9, "12-10-2020"),
    @("Haryana", 10573, -104, 130003, 1162, 1579, 7, "12-10-2020"),
    @("Himachal Pradesh", 2687, -31, 14471, 193, 250, 2, "12-10-2020"),
    @("Jammu and Kashmir", 10466, -330, 71845, 890, 1322, 9, "12-10-2020"),
    @("Jharkhand", 8167, -195, 83571, 766, 787, 3, "12-10-2020"),
    @("Karnataka", 120289, -659, 580054, 10107, 9966, 75, "12-10-2020"),
    @("Kerala", 96401, 398, 191798, 8924, 1003, 25, "12-10-2020"),
    @("Ladakh", 980, -42, 4037, 64, 64, 0, "12-10-2020"),
    @("Madhya Pradesh", 15177, -435, 129019, 1985, 2624, 25, "12-10-2020"),
    @("Maharashtra", 221637, 22, 1266240, 10461, 40349, 309, "12-10-2020"),
    @("Manipur", 2731, 123, 10504, 108, 91, 3, "12-10-2020"),
    @("Meghalaya", 2478, 41, 5142, 97, 63, 1, "12-10-2020"),
    @("Mizoram", 174, -17, 2010, 26, 0, 0, "12-10-2020"),
    @("Nagaland", 1259, 21, 5743, 49, 17, 0, "12-10-2020"),
    @("Odisha", 23602, -812, 227615, 3342, 1022, 16, "12-10-2020"),
    @("Puducherry", 4695, -24, 26291, 336, 563, 4, "12-10-2020"),
    @("Punjab", 9275, -477, 110865, 1098, 3833, 35, "12-10-2020"),
    @("Rajasthan", 21412, 58, 135990, 2072, 1650, 14, "12-10-2020"),
    @("Sikkim", 384, -66, 2920, 104, 55, 0, "12-10-2020"),
    @("Tamil Nadu", 44095, -55, 602038, 5005, 10252, 65, "12-10-2020"),
    @("Telengana", 24514, -1199, 187342, 2214, 1228, 6, "12-10-2020"),
    @("Tripura", 3742, -209, 24403, 317, 316, 1, "12-10-2020"),
    @("Uttarakhand", 7373, 52, 46931, 461, 747, 13, "12-10-2020"),
    @("Uttar Pradesh", 40019, -191, 390566, 3417, 6394, 41, "12-10-2020"),
    @("West Bengal", 30236, 443, 258948, 3110, 5622, 59, "12-10-2020")
)

$startRow = 902
for ($i = 0; $i -lt $covidData.Count; $i++) {
    $row = $covidData[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    # Leading apostrophe forces the DD-MM-YYYY text to stay a literal string
    # instead of being auto-converted to a date serial number.
    $ws.Cells.Item($r, 8).Value = "'" + $row[7]
}
